$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "J16,J15,J10,J12,J11"
